# Update workbook data ("Mise a jour du fichier via Shiny")
$wb = $excel.ActiveWorkbook

# --- Sheet "pro": update production values B2:B26 ---
$wsPro = $wb.Worksheets.Item("pro")
$proValues = @(892498.28095767635, 867204.97223946173, 843060.09939562157, 807812.33257700538, 989963.64787321072, 1004220.3389177501, 1001475.1049139029, 1016119.1603356658, 1068621.4612114637, 1048251.1507697842, 1036199.5489919075, 1074257.4916093561, 1106191.3234820557, 1180925.4847899291, 1203415.9974299127, 1023011.1413454487, 1065247.2097024948, 1163759.4496102144, 1279070.7003550082, 1343057, 1402966, 1649926.9788154145, 1518839.6399376567, 1579995.5479827698, 1700183.0408774351)
for ($i = 0; $i -lt $proValues.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $proValues[$i]
}

# --- Sheet "ind": update indicator values B2:B101 ---
$wsInd = $wb.Worksheets.Item("ind")
$indValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 774049.62066820927, 879155.68326193106, 783313.59709996183, 942880.04054559371, 962777.4721621864, 1144983.3045774491, 940900.41249388293, 951338.81076648133, 898438.59783777187, 1029680.2102070589, 1165623.0114931432, 1200017.5465738068, 1101567.3626145632, 1441398.0521946086, 1287205.6602888315, 1373635.2766218903, 1059465.1018414048, 1359038.3493043953, 1522177.7063888994, 1350365.1773913954, 1161130.9281515211, 1532582.4054490556, 1236146.7349534507, 1615877.5765377402, 1183323.0126313234, 1476526.3598635746, 2003272.8301733413, 2005040.9099544766)
for ($i = 0; $i -lt $indValues.Length; $i++) {
    $wsInd.Cells.Item($i + 2, 2).Value = $indValues[$i]
}

# --- Sheet "conso": update consumption values B2:B26 ---
$wsConso = $wb.Worksheets.Item("conso")
$consoValues = @(306828.34817878617, 298133.06067248934, 289832.18673351943, 277713.63343450398, 340334.85426003649, 345235.84711000836, 344292.69996038347, 349326.02438077383, 367375.70794501575, 360373.08678338432, 356229.99978370732, 369313.73011306248, 380292.36052552372, 405984.69533583865, 413717.31102382927, 351696.40997981455, 366216.8578556233, 400084.09016452834, 439726.68382785184, 461724, 428999, 369329.06515547639, 383404.27939192142, 398842.00977375283, 429181.22261329449)
for ($i = 0; $i -lt $consoValues.Length; $i++) {
    $wsConso.Cells.Item($i + 2, 2).Value = $consoValues[$i]
}

# "VA" sheet (B2:B26) holds =pro!B# -conso!B# formulas and recalculates
# automatically from the updated "pro" / "conso" values above.

# --- Widen column B on "pro" slightly (14.6328125 -> ~15.6328125 chars) ---
$wsPro.Columns.Item(2).ColumnWidth = 14.75

# --- Scroll "ind" sheet view down so row 66 is at the top ---
$wsInd.Activate()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1

# --- Refresh the selection rectangle on every sheet to B2:B101 ---
foreach ($name in @("pro","ind","VA","conso")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B2:B101").Select()
}

# --- Restore "pro" as the active/selected sheet (matches original state) ---
$wsPro.Activate()

Write-Output "workbook updated"
